$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '37.007.38'
Set-TextCell $ws.Range('E2') '  +1.09%  '
Set-TextCell $ws.Range('D3') '2.061.17'
Set-TextCell $ws.Range('E3') '  -1.82%  '
Set-TextCell $ws.Range('E4') '  -0.06%  '
Set-TextCell $ws.Range('D5') '249.48'
Set-TextCell $ws.Range('E5') '  -1.33%  '
Set-TextCell $ws.Range('D6') '0.672'
Set-TextCell $ws.Range('E6') '  +1.93%  '
Set-TextCell $ws.Range('E7') '  -0.05%  '
Set-TextCell $ws.Range('D8') '55.07'
Set-TextCell $ws.Range('E8') '  +11.40%  '
Set-TextCell $ws.Range('D9') '61.10'
Set-TextCell $ws.Range('E9') '  +1.12%  '
Set-TextCell $ws.Range('D10') '0.383'
Set-TextCell $ws.Range('E10') '  +1.64%  '
Set-TextCell $ws.Range('D11') '0.0803'
Set-TextCell $ws.Range('E11') '  +7.69%  '
Set-TextCell $ws.Range('E12') '  +5.83%  '
Set-TextCell $ws.Range('D13') '15.03'
Set-TextCell $ws.Range('E13') '  +2.30%  '
Set-TextCell $ws.Range('D14') '2.361.65'
Set-TextCell $ws.Range('E14') '  -1.81%  '
Set-TextCell $ws.Range('D15') '0.816'
Set-TextCell $ws.Range('E15') '  -2.27%  '
Set-TextCell $ws.Range('D16') '5.31'
Set-TextCell $ws.Range('E16') '  +3.60%  '
Set-TextCell $ws.Range('D17') '2.055.44'
Set-TextCell $ws.Range('E17') '  -2.37%  '
Set-TextCell $ws.Range('D18') '36.980.63'
Set-TextCell $ws.Range('E18') '  +1.17%  '
Set-TextCell $ws.Range('D19') '0.0₃0948'
Set-TextCell $ws.Range('E19') '  +13.50%  '
Set-TextCell $ws.Range('D20') '73.31'
Set-TextCell $ws.Range('E20') '  +0.25%  '
Set-TextCell $ws.Range('E21') '  +6.48%  '
Set-TextCell $ws.Range('D22') '5.39'
Set-TextCell $ws.Range('E22') '  +2.45%  '
Set-TextCell $ws.Range('D23') '237.34'
Set-TextCell $ws.Range('E23') '  -1.24%  '
Set-TextCell $ws.Range('E24') '  -0.04%  '
Set-TextCell $ws.Range('D25') '2.43'
Set-TextCell $ws.Range('E25') '  -4.07%  '
Set-TextCell $ws.Range('D26') '170.17'
Set-TextCell $ws.Range('E26') '  -0.56%  '
Set-TextCell $ws.Range('D27') '9.10'
Set-TextCell $ws.Range('E27') '  -1.14%  '
Set-TextCell $ws.Range('D28') '20.12'
Set-TextCell $ws.Range('E28') '  -5.10%  '
Set-TextCell $ws.Range('D29') '2.00'
Set-TextCell $ws.Range('E29') '  +0.01%  '
Set-TextCell $ws.Range('E30') '  +1.91%  '
Set-TextCell $ws.Range('E31') '  +2.32%  '
Set-TextCell $ws.Range('B32') 'ImmutableX'
Set-TextCell $ws.Range('C32') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws.Range('D32') '1.06'
Set-TextCell $ws.Range('E32') '  +8.12%  '
Set-TextCell $ws.Range('B33') 'Hedera'
Set-TextCell $ws.Range('C33') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws.Range('D33') '0.0629'
Set-TextCell $ws.Range('E33') '  +1.56%  '
Set-TextCell $ws.Range('E34') '  +6.98%  '
Set-TextCell $ws.Range('D35') '0.0892'
Set-TextCell $ws.Range('E35') '  -0.67%  '
Set-TextCell $ws.Range('D36') '0.999'
Set-TextCell $ws.Range('E36') '  -0.18%  '
Set-TextCell $ws.Range('D37') '2.27'
Set-TextCell $ws.Range('E37') '  -6.33%  '
Set-TextCell $ws.Range('E38') '  -5.02%  '
Set-TextCell $ws.Range('D39') '1.34'
Set-TextCell $ws.Range('E39') '  -0.89%  '
Set-TextCell $ws.Range('D40') '0.105'
Set-TextCell $ws.Range('E40') '  +23.69%  '
Set-TextCell $ws.Range('B41') 'VeChain'
Set-TextCell $ws.Range('C41') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws.Range('D41') '0.0225'
Set-TextCell $ws.Range('E41') '  +0.12%  '
Set-TextCell $ws.Range('B42') 'InjectiveProtocol'
Set-TextCell $ws.Range('C42') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws.Range('D42') '17.74'
Set-TextCell $ws.Range('E42') '  +7.38%  '
Set-TextCell $ws.Range('E43') '  -2.37%  '
Set-TextCell $ws.Range('D44') '96.63'
Set-TextCell $ws.Range('E44') '  -1.30%  '
Set-TextCell $ws.Range('E45') '  +0.55%  '
Set-TextCell $ws.Range('D46') '4.16'
Set-TextCell $ws.Range('E46') '  +39.69%  '
Set-TextCell $ws.Range('D47') '13.65'
Set-TextCell $ws.Range('E47') '  -52.05%  '
Set-TextCell $ws.Range('E48') '  +7.90%  '
Set-TextCell $ws.Range('D49') '1.298.32'
Set-TextCell $ws.Range('E49') '  -3.17%  '
Set-TextCell $ws.Range('E50') '  +0.99%  '
Set-TextCell $ws.Range('D51') '4.17'
Set-TextCell $ws.Range('E51') '  +8.21%  '
